$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Row 17 ("repaymentstrategy") used to hold "Mifos style" as its value;
# replace it with the new scenario label being added by this commit.
$ws.Range("B17").Value = "Penalties, Fees, Interest, Principal order"

# Adjust formatting: left/top aligned, no wrap text (distinct style from the rest of column B)
$ws.Range("B17").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B17").VerticalAlignment = -4160     # xlTop
$ws.Range("B17").WrapText = $false

# Move the selection to B17 (matches the saved selection state)
$ws.Activate()
$ws.Range("B17").Select()
